$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.09599933333334
$ws.Range("H2").Value = 60.287998
$ws.Range("I2").Value = 0.9230842418515092
$ws.Range("J2").Value = 0.9230842418515091
$ws.Range("M2").Value = 68.46613766666667
$ws.Range("N2").Value = 205.398413
$ws.Range("O2").Value = 0.4719163120948675
$ws.Range("P2").Value = 0.4719163120948675
$ws.Range("Q2").Value = 1375.895456905242
$ws.Range("R2").Value = 12383.05911214718
$ws.Range("S2").Value = 0.4356185111674509
$ws.Range("T2").Value = 0.4356185111674509
$ws.Range("G3").Value = 20.09599933333334
$ws.Range("H3").Value = 60.287998
$ws.Range("I3").Value = 0.9230842418515092
$ws.Range("J3").Value = 0.9230842418515091
$ws.Range("M3").Value = 9.278736333333333
$ws.Range("O3").Value = 0.06395551407683932
$ws.Range("P3").Value = 0.06395551407683933
$ws.Range("Q3").Value = 186.4654791688424
$ws.Range("R3").Value = 1678.189312519582
$ws.Range("S3").Value = 0.05903632722384274
$ws.Range("T3").Value = 0.05903632722384275
$ws.Range("G4").Value = 20.09599933333334
$ws.Range("H4").Value = 60.287998
$ws.Range("I4").Value = 0.9230842418515092
$ws.Range("J4").Value = 0.9230842418515091
$ws.Range("M4").Value = 67.336226
$ws.Range("N4").Value = 202.008678
$ws.Range("O4").Value = 0.4641281738282933
$ws.Range("P4").Value = 0.4641281738282933
$ws.Range("Q4").Value = 1353.188752805183
$ws.Range("R4").Value = 12178.69877524664
$ws.Range("S4").Value = 0.4284294034602156
$ws.Range("T4").Value = 0.4284294034602155
$ws.Range("I5").Value = 0.05319611498621682
$ws.Range("J5").Value = 0.05319611498621682
$ws.Range("M5").Value = 68.46613766666667
$ws.Range("N5").Value = 205.398413
$ws.Range("O5").Value = 0.4719163120948675
$ws.Range("P5").Value = 0.4719163120948675
$ws.Range("Q5").Value = 79.29102200654678
$ws.Range("R5").Value = 713.619198058921
$ws.Range("S5").Value = 0.02510411440206995
$ws.Range("T5").Value = 0.02510411440206995
$ws.Range("I6").Value = 0.05319611498621682
$ws.Range("J6").Value = 0.05319611498621682
$ws.Range("M6").Value = 9.278736333333333
$ws.Range("O6").Value = 0.06395551407683932
$ws.Range("P6").Value = 0.06395551407683933
$ws.Range("R6").Value = 96.711814144253
$ws.Range("S6").Value = 0.003402184880834153
$ws.Range("T6").Value = 0.003402184880834153
$ws.Range("I7").Value = 0.05319611498621682
$ws.Range("J7").Value = 0.05319611498621682
$ws.Range("M7").Value = 67.336226
$ws.Range("N7").Value = 202.008678
$ws.Range("O7").Value = 0.4641281738282933
$ws.Range("P7").Value = 0.4641281738282933
$ws.Range("Q7").Value = 77.98246490254733
$ws.Range("R7").Value = 701.842184122926
$ws.Range("S7").Value = 0.02468981570331272
$ws.Range("T7").Value = 0.02468981570331272
$ws.Range("G8").Value = 0.5163883333333333
$ws.Range("H8").Value = 1.549165
$ws.Range("I8").Value = 0.02371964316227407
$ws.Range("J8").Value = 0.02371964316227407
$ws.Range("M8").Value = 68.46613766666667
$ws.Range("N8").Value = 205.398413
$ws.Range("O8").Value = 0.4719163120948675
$ws.Range("P8").Value = 0.4719163120948675
$ws.Range("Q8").Value = 35.35511471946056
$ws.Range("R8").Value = 318.196032475145
$ws.Range("S8").Value = 0.01119368652534662
$ws.Range("T8").Value = 0.01119368652534662
$ws.Range("G9").Value = 0.5163883333333333
$ws.Range("H9").Value = 1.549165
$ws.Range("I9").Value = 0.02371964316227407
$ws.Range("J9").Value = 0.02371964316227407
$ws.Range("M9").Value = 9.278736333333333
$ws.Range("O9").Value = 0.06395551407683932
$ws.Range("P9").Value = 0.06395551407683933
$ws.Range("Q9").Value = 4.791431190609444
$ws.Range("R9").Value = 43.122880715485
$ws.Range("S9").Value = 0.001517001972162425
$ws.Range("T9").Value = 0.001517001972162425
$ws.Range("G10").Value = 0.5163883333333333
$ws.Range("H10").Value = 1.549165
$ws.Range("I10").Value = 0.02371964316227407
$ws.Range("J10").Value = 0.02371964316227407
$ws.Range("M10").Value = 67.336226
$ws.Range("N10").Value = 202.008678
$ws.Range("O10").Value = 0.4641281738282933
$ws.Range("P10").Value = 0.4641281738282933
$ws.Range("Q10").Value = 34.77164151709667
$ws.Range("R10").Value = 312.94477365387
$ws.Range("S10").Value = 0.01100895466476503
$ws.Range("T10").Value = 0.01100895466476503
